$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update analysis values for D2893 sample rows (reprocessed results).
# Each row corresponds to one analysis (rows 2-22); only the cells whose
# underlying values changed due to reprocessing are updated below.

# Row 2
$ws.Range("E2").Value = 100.07
$ws.Range("H2").Value = 0.4
$ws.Range("I2").Value = 100.07
$ws.Range("J2").Value = 0.34
$ws.Range("K2").Value = 15.61
$ws.Range("L2").Value = 0.4
$ws.Range("M2").Value = 16.02

# Row 3
$ws.Range("E3").Value = 100.07
$ws.Range("I3").Value = 100.07
$ws.Range("K3").Value = 15.67
$ws.Range("M3").Value = 16.08

# Row 4
$ws.Range("E4").Value = 100.07
$ws.Range("I4").Value = 100.07
$ws.Range("K4").Value = 15.44
$ws.Range("M4").Value = 15.86

# Row 5
$ws.Range("E5").Value = 100.07
$ws.Range("H5").Value = 0.36
$ws.Range("I5").Value = 100.07
$ws.Range("J5").Value = 0.3
$ws.Range("K5").Value = 20.09
$ws.Range("L5").Value = 0.36
$ws.Range("M5").Value = 20.42

# Row 6
$ws.Range("E6").Value = 100.07
$ws.Range("I6").Value = 100.07
$ws.Range("K6").Value = 15.99
$ws.Range("M6").Value = 16.4

# Row 7
$ws.Range("E7").Value = 100.07
$ws.Range("H7").Value = 0.21
$ws.Range("I7").Value = 100.07
$ws.Range("K7").Value = 59.51
$ws.Range("L7").Value = 0.21
$ws.Range("M7").Value = 59.62

# Row 8
$ws.Range("E8").Value = 100.07
$ws.Range("H8").Value = 0.42
$ws.Range("I8").Value = 100.07
$ws.Range("J8").Value = 0.35
$ws.Range("K8").Value = 18.13
$ws.Range("L8").Value = 0.42
$ws.Range("M8").Value = 18.49

# Row 9
$ws.Range("E9").Value = 100.07
$ws.Range("H9").Value = 0.46
$ws.Range("I9").Value = 100.07
$ws.Range("J9").Value = 0.38
$ws.Range("K9").Value = 16.03
$ws.Range("L9").Value = 0.46
$ws.Range("M9").Value = 16.43

# Row 10
$ws.Range("E10").Value = 100.07
$ws.Range("H10").Value = 0.31
$ws.Range("I10").Value = 100.07
$ws.Range("J10").Value = 0.26
$ws.Range("K10").Value = 20.08
$ws.Range("L10").Value = 0.31
$ws.Range("M10").Value = 20.4

# Row 11
$ws.Range("E11").Value = 100.07
$ws.Range("H11").Value = 0.29
$ws.Range("I11").Value = 100.07
$ws.Range("J11").Value = 0.24
$ws.Range("K11").Value = 25.62
$ws.Range("L11").Value = 0.29
$ws.Range("M11").Value = 25.87

# Row 12
$ws.Range("E12").Value = 100.07
$ws.Range("H12").Value = 0.44
$ws.Range("I12").Value = 100.07
$ws.Range("J12").Value = 0.37
$ws.Range("K12").Value = 18.79
$ws.Range("L12").Value = 0.44
$ws.Range("M12").Value = 19.14

# Row 13
$ws.Range("E13").Value = 100.07
$ws.Range("I13").Value = 100.07
$ws.Range("K13").Value = 25.07
$ws.Range("M13").Value = 25.33

# Row 14
$ws.Range("E14").Value = 100.07
$ws.Range("H14").Value = 0.42
$ws.Range("I14").Value = 100.07
$ws.Range("J14").Value = 0.35
$ws.Range("K14").Value = 15.66
$ws.Range("L14").Value = 0.42
$ws.Range("M14").Value = 16.07

# Row 15
$ws.Range("E15").Value = 100.07
$ws.Range("H15").Value = 0.41
$ws.Range("I15").Value = 100.07
$ws.Range("J15").Value = 0.34
$ws.Range("K15").Value = 16.82
$ws.Range("L15").Value = 0.41
$ws.Range("M15").Value = 17.21

# Row 16
$ws.Range("E16").Value = 100.07
$ws.Range("I16").Value = 100.07
$ws.Range("K16").Value = 18.51
$ws.Range("M16").Value = 18.86

# Row 17
$ws.Range("E17").Value = 100.07
$ws.Range("H17").Value = 0.35
$ws.Range("I17").Value = 100.07
$ws.Range("J17").Value = 0.29
$ws.Range("K17").Value = 21.5
$ws.Range("L17").Value = 0.35
$ws.Range("M17").Value = 21.8

# Row 18
$ws.Range("E18").Value = 100.07
$ws.Range("I18").Value = 100.07
$ws.Range("K18").Value = 49.17
$ws.Range("M18").Value = 49.3

# Row 19
$ws.Range("E19").Value = 100.07
$ws.Range("H19").Value = 0.15
$ws.Range("I19").Value = 100.07
$ws.Range("J19").Value = 0.12
$ws.Range("K19").Value = 96.43000000000001
$ws.Range("L19").Value = 0.15
$ws.Range("M19").Value = 96.5

# Row 20
$ws.Range("E20").Value = 100.07
$ws.Range("H20").Value = 0.25
$ws.Range("I20").Value = 100.07
$ws.Range("J20").Value = 0.2
$ws.Range("K20").Value = 54.64
$ws.Range("L20").Value = 0.25
$ws.Range("M20").Value = 54.76

# Row 21
$ws.Range("E21").Value = 100.07
$ws.Range("I21").Value = 100.07
$ws.Range("K21").Value = 86.81
$ws.Range("M21").Value = 86.88

# Row 22
$ws.Range("E22").Value = 100.07
$ws.Range("H22").Value = 0.36
$ws.Range("I22").Value = 100.07
$ws.Range("J22").Value = 0.3
$ws.Range("K22").Value = 37.72
$ws.Range("L22").Value = 0.36
$ws.Range("M22").Value = 37.89
